# Rewrites the monthly PPI data rows (rows 2-49) so each year block
# starts with Oct/Nov/Dec of that year followed by Jan-Sep, matching the
# re-ordering applied in the target workbook. Row 1 (headers) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @(ColA_label, ColB_value, ColC_value, ColD_value)
$data = @{
    2 = @("2014-10", 106.725, 101.5483, 98.559)
    3 = @("2014-11", 106.6044, 101.1426, 100.6241)
    4 = @("2014-12", 106.0245, 101.4391, 98.2025)
    5 = @("2014-01", 106.7637, 100.8731, 96.0581)
    6 = @("2014-02", 107.2275, 101.0438, 96.6465)
    7 = @("2014-03", 107.0932, 100.9209, 97.6887)
    8 = @("2014-04", 107.4918, 100.9206, 96.7529)
    9 = @("2014-05", 109.0628, 100.986, 94.4473)
    10 = @("2014-06", 108.719, 100.5403, 96.6622)
    11 = @("2014-07", 109.9193, 101.0631, 99.33969999999999)
    12 = @("2014-08", 109.7716, 101.4661, 99.7051)
    13 = @("2014-09", 107.659, 101.8374, 99.9813)
    14 = @("2015-10", 100.5, 100.7, 78.3)
    15 = @("2015-11", 100.7158, 100.6953, 77.2026)
    16 = @("2015-12", 100.9756, 100.6681, 77.26309999999999)
    17 = @("2015-01", 104.7694, 101.0813, 79.27460000000001)
    18 = @("2015-02", 103.8162, 101.133, 78.027)
    19 = @("2015-03", 103.8615, 100.903, 76.0416)
    20 = @("2015-04", 102.4568, 101.1407, 83.22799999999999)
    21 = @("2015-05", 100.8059, 101.1071, 80.1913)
    22 = @("2015-06", 101.2152, 101.5267, 79.82859999999999)
    23 = @("2015-07", 99.4438, 100.9969, 78.1323)
    24 = @("2015-08", 99.4987, 100.5067, 77.6249)
    25 = @("2015-09", 100.4316, 100.3657, 75.3171)
    26 = @("2016-10", 98.59999999999999, 101.7, 97.09999999999999)
    27 = @("2016-11", 98.40000000000001, 102.3, 97.90000000000001)
    28 = @("2016-12", 99.09999999999999, 102.8, 102.4)
    29 = @("2016-01", 101.436, 100.9369, 96.265)
    30 = @("2016-02", 102.0294, 100.71, 95.98569999999999)
    31 = @("2016-03", 101.5128, 100.7751, 95.40219999999999)
    32 = @("2016-04", 101.3024, 99.845, 94.0009)
    33 = @("2016-05", 101.8, 100.2, 94.8)
    34 = @("2016-06", 101.1, 100.7, 96.2)
    35 = @("2016-07", 101.2, 101.2, 96.2)
    36 = @("2016-08", 101.1, 101.6, 96.3)
    37 = @("2016-09", 99, 101.4, 96.59999999999999)
    38 = @("2017-10", 100, 103.5, 102.9)
    39 = @("2017-11", 99.8, 103.3, 103.6)
    40 = @("2017-12", 99.59999999999999, 103, 100.1)
    41 = @("2017-01", 98.5, 103.5, 104.6)
    42 = @("2017-02", 98.7, 103.5, 104.9)
    43 = @("2017-03", 99, 104.3, 107)
    44 = @("2017-04", 99.5, 104.9, 106.6)
    45 = @("2017-05", 99.2, 105, 105.6)
    46 = @("2017-06", 99.09999999999999, 104.2, 100)
    47 = @("2017-07", 98.59999999999999, 104, 100)
    48 = @("2017-08", 98.8, 103.7, 102.3)
    49 = @("2017-09", 99.8, 103.8, 102.1)
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
